# Logged Week 15 and simulated Week 16
# Appends the new per-game log numbers to the running text logs on the
# YDS and ST sheets, and bumps the season-total numeric cells on
# OFF / DEF / ST / TURNS / PEN to reflect the two additional games.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# YDS sheet: append newly logged per-play yardage numbers
# ---------------------------------------------------------------------
$ydsWs = $wb.Worksheets("YDS")

$ydsWs.Range("B2").Value = $ydsWs.Range("B2").Value2 + " 10 0 7 6 -3 5 5 0 27 6 8 0 0 -3 1 0 1 4 5 3 0 6 0 2 8 1 16 5 0 2 1 2 -2 3"
$ydsWs.Range("B3").Value = $ydsWs.Range("B3").Value2 + " 4 5 15 5 13 5 13 5 37 5 15 14 8 22 16 6 5 8 11 -2"
$ydsWs.Range("C2").Value = $ydsWs.Range("C2").Value2 + " 0 6 2 5 1 22 0 3 23 6 1 2 1 6 3 0 4 12"
$ydsWs.Range("C3").Value = $ydsWs.Range("C3").Value2 + " 6 -1 14 13 2 5 24 10 15 0 3 9 20 19 27 16 5 7 6 7 13 11 26"

# ---------------------------------------------------------------------
# OFF sheet: season totals now include the two new games
# ---------------------------------------------------------------------
$offWs = $wb.Worksheets("OFF")

$offWs.Range("B2").Value = 5
$offWs.Range("C2").Value = 330
$offWs.Range("E2").Value = 21
$offWs.Range("F2").Value = 114
$offWs.Range("G2").Value = 92
$offWs.Range("I2").Value = 12
$offWs.Range("J2").Value = 42
$offWs.Range("L2").Value = 544
$offWs.Range("M2").Value = 363
$offWs.Range("O2").Value = 51
$offWs.Range("P2").Value = 24
$offWs.Range("Q2").Value = 920

$offWs.Range("C3").Value = 346
$offWs.Range("E3").Value = 47
$offWs.Range("F3").Value = 178
$offWs.Range("G3").Value = 82
$offWs.Range("H3").Value = 50
$offWs.Range("I3").Value = 114
$offWs.Range("J3").Value = 93
$offWs.Range("N3").Value = 40

# ---------------------------------------------------------------------
# DEF sheet: season totals now include the two new games
# ---------------------------------------------------------------------
$defWs = $wb.Worksheets("DEF")

$defWs.Range("C2").Value = 439
$defWs.Range("F2").Value = 128
$defWs.Range("G2").Value = 114
$defWs.Range("J2").Value = 61
$defWs.Range("L2").Value = 509
$defWs.Range("M2").Value = 341
$defWs.Range("O2").Value = 27
$defWs.Range("Q2").Value = 959

$defWs.Range("B3").Value = 20
$defWs.Range("C3").Value = 289
$defWs.Range("E3").Value = 50
$defWs.Range("F3").Value = 190
$defWs.Range("G3").Value = 78
$defWs.Range("H3").Value = 55
$defWs.Range("I3").Value = 91
$defWs.Range("J3").Value = 90
$defWs.Range("N3").Value = 27

# ---------------------------------------------------------------------
# ST sheet: kicking totals plus the logged per-game special-teams data
# ---------------------------------------------------------------------
$stWs = $wb.Worksheets("ST")

$stWs.Range("B2").Value = 127
$stWs.Range("D2").Value = 115
$stWs.Range("F2").Value = 126
$stWs.Range("G2").Value = 118
$stWs.Range("J2").Value = 56
$stWs.Range("K2").Value = 52
$stWs.Range("L2").Value = 27
$stWs.Range("M2").Value = 18

$stWs.Range("B3").Value = 50

$stWs.Range("D3").Value = $stWs.Range("D3").Value2 + " 40 55 38"
$stWs.Range("B4").Value = $stWs.Range("B4").Value2 + " 65 65"
$stWs.Range("D4").Value = $stWs.Range("D4").Value2 + " 6 12 0"
$stWs.Range("B5").Value = $stWs.Range("B5").Value2 + " 17 31"
$stWs.Range("D5").Value = $stWs.Range("D5").Value2 + " 0 10 14"
$stWs.Range("B6").Value = $stWs.Range("B6").Value2 + " 17"

# ---------------------------------------------------------------------
# TURNS sheet: turnover totals
# ---------------------------------------------------------------------
$turnsWs = $wb.Worksheets("TURNS")

$turnsWs.Range("C2").Value = 5
$turnsWs.Range("D2").Value = 12
$turnsWs.Range("E2").Value = 14
$turnsWs.Range("D3").Value = 11

# ---------------------------------------------------------------------
# PEN sheet: penalty totals
# ---------------------------------------------------------------------
$penWs = $wb.Worksheets("PEN")

$penWs.Range("B2").Value = 35
$penWs.Range("B3").Value = 28
$penWs.Range("D4").Value = 22
